$d = $word.ActiveDocument

$replacements = @(
    @("2024-06-17 Monday", "2024-06-18 Tuesday"),
    @("517÷2=258, 1", "797÷4=199, 1"),
    @("735÷7=105, 0", "894÷8=111, 6"),
    @("901÷7=128, 5", "978÷9=108, 6"),
    @("961÷4=240, 1", "459÷6=76, 3"),
    @("562÷3=187, 1", "190÷8=23, 6"),
    @("994÷9=110, 4", "396÷2=198, 0"),
    @("978÷6=163, 0", "107÷2=53, 1"),
    @("141÷9=15, 6", "320÷9=35, 5"),
    @("826÷5=165, 1", "980÷4=245, 0"),
    @("273÷4=68, 1", "526÷5=105, 1"),
    @("870÷9=96, 6", "622÷3=207, 1"),
    @("139÷9=15, 4", "985÷3=328, 1"),
    @("127÷2=63, 1", "319÷4=79, 3"),
    @("165÷5=33, 0", "800÷6=133, 2"),
    @("582÷6=97, 0", "816÷7=116, 4"),
    @("306÷8=38, 2", "377÷3=125, 2"),
    @("183÷6=30, 3", "456÷8=57, 0"),
    @("666÷5=133, 1", "410÷3=136, 2"),
    @("919÷3=306, 1", "626÷4=156, 2"),
    @("856÷5=171, 1", "185÷8=23, 1"),
    @("645÷6=107, 3", "848÷7=121, 1"),
    @("752÷9=83, 5", "179÷6=29, 5"),
    @("542÷3=180, 2", "386÷9=42, 8"),
    @("919÷8=114, 7", "732÷4=183, 0"),
    @("892÷6=148, 4", "884÷6=147, 2"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
